$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (Förändrad) holds a date serial that moved from 2023-10-04 (45203)
# to 2023-10-06 (45205) for every data row (rows 2 through 61).
for ($row = 2; $row -le 61; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45205
    }
}
